$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Thu Feb  9 08:50:17 2023"

$ws.Range("A4").Value = "29P188276"
$ws.Range("B4").Value = "Thu Feb  9 09:13:19 2023"
$ws.Range("C4").Value = "Thu Feb  9 09:13:19 2023"

$ws.Range("A5").Value = "29P188276"
$ws.Range("B5").Value = "Thu Feb  9 09:13:43 2023"
$ws.Range("C5").Value = "Thu Feb  9 09:13:43 2023"

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "123456789"
$ws.Range("B6").Value = "Fri Feb 10 09:56:40 2023"
$ws.Range("C6").Value = "Fri Feb 10 09:56:40 2023"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "123456789"
$ws.Range("B7").Value = "Fri Feb 10 10:04:28 2023"
$ws.Range("C7").Value = "Fri Feb 10 10:04:28 2023"
